$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that change for each data row (row index is column A, values 0..23)
$cols = @("B","C","D","E","F","J","L","N","O")

# Each entry: worksheet row number, followed by new values for columns B,C,D,E,F,J,L,N,O
$data = @(
  ,@(2, 1.606891053220522, 0.1429500831398798, 0.1249724091644211, 0.120966038716622, 1.675643970990265, 0.148110158128766, 0.3160357423338951, 1.454331322397422, 4.329673273218816)
  ,@(3, 1.510906854829216, 0.1320061420179144, 0.1243874862335375, 0.1215383176540445, 1.678269952518946, 0.1493628627989612, 0.3092004679035654, 1.467457292950066, 4.339556449032727)
  ,@(4, 1.45240851255312, 0.1252392689780493, 0.1240625125239561, 0.121927452672498, 1.680778051694631, 0.1501807542709148, 0.3051170192270689, 1.476064786199359, 4.348080734095475)
  ,@(5, 1.42868126370513, 0.1224699808846736, 0.1239387148678404, 0.122095538650683, 1.682025394674994, 0.1505263150356786, 0.3034816613593705, 1.479710325444593, 4.352171853911273)
  ,@(6, 1.424748140946463, 0.1220094384884618, 0.1239186808511406, 0.12212402409512, 1.6822461218947, 0.1505844361339594, 0.3032118482649651, 1.480323997060786, 4.35288846605286)
  ,@(7, 1.452088065832527, 0.1252019686967714, 0.1240608079437564, 0.1219296810094566, 1.68079396163904, 0.1501853649555631, 0.3050948478849023, 1.476113392675984, 4.34813340875408)
  ,@(8, 1.573706052430737, 0.1391865047557275, 0.1247636585375034, 0.1211555348524058, 1.676363514665823, 0.148531984280444, 0.3136554864498606, 1.458743455477808, 4.33257114206728)
  ,@(9, 1.815604218637304, 0.1662299329418602, 0.1264115142573488, 0.1199362711098431, 1.674783178183148, 0.1456758152408106, 0.3313372079276604, 1.429025804156126, 4.321553114353719)
  ,@(10, 1.995345918230043, 0.1858616385550818, 0.1277846476447806, 0.1192217142890986, 1.677958678469295, 0.1438120303167043, 0.3448670973744328, 1.409835277288408, 4.325369774305301)
  ,@(11, 2.077542605579652, 0.1947401106865527, 0.1284442136089794, 0.1189358017912401, 1.680345898222058, 0.1430149403840879, 0.3511380776938466, 1.401677819984883, 4.329698269913365)
  ,@(12, 2.108729037880778, 0.1980945399575944, 0.128698959435944, 0.1188331467689512, 1.681385461308736, 0.1427203910827739, 0.3535293053748916, 1.398671050879194, 4.331710474109173)
  ,@(13, 2.102009825291191, 0.1973724466961926, 0.1286438743856024, 0.1188550059240452, 1.681155542750005, 0.142783503373014, 0.3530135784904189, 1.399314953254709, 4.331260511081894)
  ,@(14, 2.080107132670094, 0.1950162359719343, 0.1284650720677476, 0.1189272438760565, 1.680428706466458, 0.142990561546906, 0.3513344747977953, 1.40142880235625, 4.329856335853123)
  ,@(15, 2.06669890783121, 0.1935719870139962, 0.1283561981984604, 0.1189722223544454, 1.680001154694608, 0.1431183399958993, 0.350308124989553, 1.402734310861774, 4.32904483601763)
  ,@(16, 1.98998270561674, 0.1852803472310711, 0.1277422428105197, 0.119241185726624, 1.677821638668334, 0.1438651427261597, 0.3444595978460256, 1.410379902070538, 4.325139088433104)
  ,@(17, 1.943029061000516, 0.1801802358880025, 0.1273745172685707, 0.1194162008901039, 1.676726030552089, 0.1443362764876586, 0.3409013540575643, 1.415216810345875, 4.323407218344244)
  ,@(18, 1.91606326184592, 0.1772418936928659, 0.1271663004902308, 0.119520550191103, 1.676184592959885, 0.1446120379055031, 0.3388656900113176, 1.418052751354907, 4.322655070123545)
  ,@(19, 1.906940146557133, 0.1762461861647182, 0.127096367987015, 0.1195565144639197, 1.676016509872767, 0.1447062267852903, 0.3381783336393767, 1.419022207039468, 4.322442299007662)
  ,@(20, 1.948023157047203, 0.1807236590578043, 0.1274133221408107, 0.1193971889300869, 1.676833476824228, 0.1442856290687331, 0.3412790037223061, 1.41469633711975, 4.323566325699801)
  ,@(21, 2.086538860493306, 0.1957085212566199, 0.1285174557292947, 0.1189058735988251, 1.680638516279657, 0.1429295457132707, 0.3518272205487136, 1.400805680704728, 4.330258647778379)
  ,@(22, 2.177417753308873, 0.2054572969296089, 0.1292680947186042, 0.1186174855078832, 1.683915572558973, 0.1420857632086818, 0.3588174439683911, 1.39220693177986, 4.336807314771761)
  ,@(23, 2.128882429502198, 0.2002583361700658, 0.1288648217416295, 0.118768415141238, 1.682094231954196, 0.1425322197322973, 0.3550778668585224, 1.396752370912694, 4.333113059958919)
  ,@(24, 1.945765237390901, 0.180477996804882, 0.1273957685085918, 0.1194057726153659, 1.676784624871644, 0.1443085115019738, 0.3411082370938487, 1.414931471307476, 4.323493634707035)
  ,@(25, 1.749805672564207, 0.1589552241981096, 0.125937061626658, 0.1202342195758916, 1.674449407329206, 0.1464072261235145, 0.3264587194025381, 1.43660068455479, 4.322443730337767)
)

foreach ($entry in $data) {
  $row = $entry[0]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $cellRef = "{0}{1}" -f $cols[$i], $row
    $ws.Range($cellRef).Value = $entry[$i + 1]
  }
}

Write-Host "Updated" $data.Length "rows"